$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so values like "242.02" or
# "0.9997" are not silently reinterpreted as numbers (losing formatting,
# trailing zeros, etc.), matching the inline-string cell type used in the file.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "29.480.11"
$ws.Range("E2").Value = "  -1.11%  "

# Row 3
$ws.Range("D3").Value = "1.849.26"
$ws.Range("E3").Value = "  -0.64%  "

# Row 4
$ws.Range("E4").Value = "  -0.05%  "

# Row 5
$ws.Range("D5").Value = "242.02"
$ws.Range("E5").Value = "  -1.34%  "

# Row 6
$ws.Range("D6").Value = "0.6275"
$ws.Range("E6").Value = "  -2.04%  "

# Row 7
$ws.Range("D7").Value = "0.9997"
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").Value = "48.23"
$ws.Range("E8").Value = "  +1.49%  "

# Row 9
$ws.Range("D9").Value = "0.07545"
$ws.Range("E9").Value = "  -0.01%  "

# Row 10
$ws.Range("D10").Value = "0.2978"
$ws.Range("E10").Value = "  -0.06%  "

# Row 11
$ws.Range("E11").Value = "  -1.35%  "

# Row 12
$ws.Range("D12").Value = "1.941.06"
$ws.Range("E12").Value = "  +4.37%  "

# Row 13
$ws.Range("D13").Value = "0.07698"
$ws.Range("E13").Value = "  +0.33%  "

# Row 14
$ws.Range("D14").Value = "5.006"
$ws.Range("E14").Value = "  -0.86%  "

# Row 15
$ws.Range("D15").Value = "0.6857"
$ws.Range("E15").Value = "  -1.05%  "

# Row 16
$ws.Range("D16").Value = "83.84"
$ws.Range("E16").Value = "  -0.36%  "

# Row 17
$ws.Range("D17").Value = "0.000009770"
$ws.Range("E17").Value = "  -0.26%  "

# Row 18
$ws.Range("D18").Value = "2.181.46"
$ws.Range("E18").Value = "  +3.41%  "

# Row 19
$ws.Range("D19").Value = "6.214"
$ws.Range("E19").Value = "  +1.49%  "

# Row 20
$ws.Range("D20").Value = "29.576.73"
$ws.Range("E20").Value = "  -0.80%  "

# Row 21
$ws.Range("D21").Value = "234.41"
$ws.Range("E21").Value = "  -0.97%  "

# Row 22
$ws.Range("D22").Value = "12.49"
$ws.Range("E22").Value = "  -1.45%  "

# Row 23
$ws.Range("D23").Value = "0.9994"

# Row 24
$ws.Range("D24").Value = "7.634"
$ws.Range("E24").Value = "  +1.63%  "

# Row 25
$ws.Range("D25").Value = "0.9999"
$ws.Range("E25").Value = "  -0.07%  "

# Row 26
$ws.Range("D26").Value = "155.38"
$ws.Range("E26").Value = "  -2.46%  "

# Row 27
$ws.Range("E27").Value = "  -1.98%  "

# Row 28
$ws.Range("D28").Value = "8.430"
$ws.Range("E28").Value = "  -1.61%  "

# Row 29
$ws.Range("E29").Value = "  -1.38%  "

# Row 30
$ws.Range("D30").Value = "1.478"
$ws.Range("E30").Value = "  -1.51%  "

# Row 31
$ws.Range("D31").Value = "0.05852"
$ws.Range("E31").Value = "  -5.34%  "

# Row 32
$ws.Range("D32").Value = "1.265"
$ws.Range("E32").Value = "  -2.53%  "

# Row 33
$ws.Range("E33").Value = "  -1.48%  "

# Row 34
$ws.Range("D34").Value = "4.023"
$ws.Range("E34").Value = "  -2.21%  "

# Row 35
$ws.Range("D35").Value = "1.883"
$ws.Range("E35").Value = "  -0.91%  "

# Row 36
$ws.Range("D36").Value = "1.170"
$ws.Range("E36").Value = "  -0.59%  "

# Row 37
$ws.Range("D37").Value = "0.7184"
$ws.Range("E37").Value = "  -1.88%  "

# Row 38
$ws.Range("E38").Value = "  -0.83%  "

# Row 39
$ws.Range("E39").Value = "  -0.52%  "

# Row 40
$ws.Range("D40").Value = "1.235.39"
$ws.Range("E40").Value = "  +1.90%  "

# Row 41
$ws.Range("D41").Value = "0.01776"
$ws.Range("E41").Value = "  -0.60%  "

# Row 42
$ws.Range("D42").Value = "0.9081"
$ws.Range("E42").Value = "  -1.62%  "

# Row 43
$ws.Range("D43").Value = "6.129"
$ws.Range("E43").Value = "  -2.62%  "

# Row 44
$ws.Range("D44").Value = "2.089.61"
$ws.Range("E44").Value = "  +3.45%  "

# Row 45
$ws.Range("D45").Value = "0.9995"
$ws.Range("E45").Value = "  -0.07%  "

# Row 46
$ws.Range("D46").Value = "101.94"
$ws.Range("E46").Value = "  -0.17%  "

# Row 47
$ws.Range("D47").Value = "67.33"
$ws.Range("E47").Value = "  +0.87%  "

# Row 48
$ws.Range("D48").Value = "7.304"
$ws.Range("E48").Value = "  +8.64%  "

# Row 49
$ws.Range("B49").Value = "TheSandbox"
$ws.Range("C49").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D49").Value = "0.4028"
$ws.Range("E49").Value = "  -0.86%  "

# Row 50
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "9.131"
$ws.Range("E50").Value = "  -1.05%  "

# Row 51
$ws.Range("D51").Value = "1.711"
$ws.Range("E51").Value = "  +2.46%  "
